$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.275.21"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "3.297.67"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "556.76"
$ws.Range("E5").Value = "  -3.43%  "
$ws.Range("D6").Value = "141.86"
$ws.Range("E6").Value = "  -4.97%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.296.65"
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("D10").Value = "7.83"
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "3.862.88"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "26.79"
$ws.Range("E15").Value = "  -5.69%  "
$ws.Range("D16").Value = "3.300.68"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "60.301.04"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("D20").Value = "13.98"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").Value = "8.62"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").Value = "374.69"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.40"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "0.534"
$ws.Range("E25").Value = "  -5.46%  "
$ws.Range("D26").Value = "3.442.47"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("E27").Value = "  -8.74%  "
$ws.Range("E28").Value = "  -5.04%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "7.15"
$ws.Range("E30").Value = "  -6.46%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("D33").Value = "7.56"
$ws.Range("E33").Value = "  -4.39%  "
$ws.Range("D34").Value = "22.56"
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  -7.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.10"
$ws.Range("E36").Value = "  -6.48%  "
$ws.Range("D37").Value = "166.59"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").Value = "6.67"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").Value = "26.84"
$ws.Range("E40").Value = "  -13.73%  "
$ws.Range("D41").Value = "3.328.72"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "0.0726"
$ws.Range("E42").Value = "  -6.63%  "
$ws.Range("D43").Value = "41.91"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.750"
$ws.Range("E44").Value = "  -3.37%  "
$ws.Range("D45").Value = "4.12"
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").Value = "1.57"
$ws.Range("E46").Value = "  -5.91%  "
$ws.Range("E47").Value = "  -4.39%  "
$ws.Range("D48").Value = "2.352.25"
$ws.Range("E48").Value = "  -7.59%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "6.39"
$ws.Range("E50").Value = "  -6.82%  "
$ws.Range("D51").Value = "21.27"
$ws.Range("E51").Value = "  -5.19%  "
